# Update TPM-derived values in the Wnt2-Fzd9 LR-pair sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 0.976531
$ws.Range("N2").Value = 1.953062
$ws.Range("O2").Value = 0.3681359341666424
$ws.Range("P2").Value = 0.3330144816631303
$ws.Range("Q2").Value = 0.561054167678
$ws.Range("R2").Value = 3.366325006068
$ws.Range("S2").Value = 0.3681359341666424
$ws.Range("T2").Value = 0.3330144816631303

# Row 3
$ws.Range("O3").Value = 0.09773758465004194
$ws.Range("P3").Value = 0.1326196171188222
$ws.Range("S3").Value = 0.09773758465004194
$ws.Range("T3").Value = 0.1326196171188222

# Row 4
$ws.Range("M4").Value = 0.127556
$ws.Range("N4").Value = 0.382668
$ws.Range("O4").Value = 0.04808648902959583
$ws.Range("P4").Value = 0.06524830531189832
$ws.Range("Q4").Value = 0.073285769128
$ws.Range("R4").Value = 0.659571922152
$ws.Range("S4").Value = 0.04808648902959583
$ws.Range("T4").Value = 0.06524830531189832

# Row 5
$ws.Range("M5").Value = 1.116584
$ws.Range("N5").Value = 2.233168
$ws.Range("O5").Value = 0.420933584203191
$ws.Range("P5").Value = 0.380775051681252
$ws.Range("Q5").Value = 0.641519938192
$ws.Range("R5").Value = 3.849119629152
$ws.Range("S5").Value = 0.420933584203191
$ws.Range("T5").Value = 0.380775051681252

# Row 6
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.03536566666666666
$ws.Range("N6").Value = 0.106097
$ws.Range("O6").Value = 0.01333226772704545
$ws.Range("P6").Value = 0.01809048430670052
$ws.Range("Q6").Value = 0.02031891939533333
$ws.Range("R6").Value = 0.182870274558
$ws.Range("S6").Value = 0.01333226772704545
$ws.Range("T6").Value = 0.01809048430670052

# Row 7
$ws.Range("M7").Value = 0.137338
$ws.Range("N7").Value = 0.412014
$ws.Range("O7").Value = 0.05177414022348326
$ws.Range("P7").Value = 0.07025205991819664
$ws.Range("Q7").Value = 0.078905899844
$ws.Range("R7").Value = 0.710153098596
$ws.Range("S7").Value = 0.05177414022348326
$ws.Range("T7").Value = 0.07025205991819664
